$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update existing values in column C (Número de cadernos)
$ws.Range("C2").Value = 449
$ws.Range("C4").Value = 2037
$ws.Range("C7").Value = 2353
$ws.Range("C8").Value = 2760
$ws.Range("C9").Value = 3030
$ws.Range("C14").Value = 4743
$ws.Range("C15").Value = 9175
$ws.Range("C17").Value = 1367
$ws.Range("C18").Value = 2057
$ws.Range("C19").Value = 1555
$ws.Range("C20").Value = 856

# Add new row 21 - copy formatting from row 20, then set values
$ws.Range("A20:C20").Copy()
$ws.Range("A21:C21").PasteSpecial(-4122)  # xlPasteFormats

$ws.Cells.Item(21, 1).Value = 19
$ws.Cells.Item(21, 2).Value = 44409
$ws.Cells.Item(21, 3).Value = 25
